$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEAVE CREDITS")
$ws.Rows(20).Insert()
